$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for the Vietnamese header labels (shifts all
# existing data rows down by one, same as the author clicking
# Insert > Insert Sheet Rows on row 1 and then typing the headers).
$ws.Rows.Item(1).Insert()

$headers = @("Mã sản phẩm", "Tên sản phẩm", "Hình ảnh", "CPU", "RAM", "ROM", "Card đồ hoạ", "Màn hình", "Pin", "Hãng", "Giá", "Tình trạng")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The new header text widened column A ("Mã sản phẩm") and column K ("Giá")
# beyond their previous (unset) widths, so auto-fit those columns like Excel
# would after the import.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(11).AutoFit() | Out-Null

# Restore the cursor/viewport to where the author ended up after the edit.
$ws.Range("I20").Select() | Out-Null
